# Fix autofilling issue: replace the old Windows username ("AL-Thuraya")
# baked into the absolute sample-path formulas with the new one ("moham"),
# across every formula/cell on the active sheet, then restore the
# last-used selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the stale username in every formula (and any literal text) on
# the sheet's used range - this updates both the formula source and the
# cached/evaluated value for cells such as B2 and B3.
$used = $ws.UsedRange
$used.Replace("AL-Thuraya", "moham", -4163, 1, $false, $false, $true) | Out-Null

# Update the active selection to match the author's last cursor position.
$ws.Range("B6").Select() | Out-Null
